$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings; runs share identical
# formatting so a plain replacement of the full visible text reproduces the
# same rendered content) ---
$ws.Range("A8").Value = "Volume 31   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  7/1/2024  Through  7/7/2024"

# --- Cells that flip between "blank placeholder" (text " 0 " shared string)
# and a real number need their style swapped too, not just the value.
# Copying from a same-column-kind cell brings the right style across, then
# we overwrite the value as needed.

# C16: was the text placeholder -> becomes the number 3 (style 14 -> 15)
$ws.Range("D16").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 3

# C22: was the number 1 -> becomes the text placeholder (style 15 -> 14)
$ws.Range("D15").Copy($ws.Range("C22"))

# C28: was the number 3 -> becomes the text placeholder (style 15 -> 14)
$ws.Range("D15").Copy($ws.Range("C28"))

# --- Remaining plain numeric value updates ---
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 4
$ws.Range("I15").Value = 8
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = -20
$ws.Range("M15").Value = 60
$ws.Range("N15").Value = -33.333333333333

$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 16
$ws.Range("H16").Value = 77.777777777777
$ws.Range("I16").Value = 71
$ws.Range("J16").Value = 72
$ws.Range("K16").Value = -1.388888888888
$ws.Range("L16").Value = 16.393442622950
$ws.Range("M16").Value = -14.457831325301
$ws.Range("N16").Value = -82.555282555282

$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 46.666666666666
$ws.Range("I17").Value = 137
$ws.Range("J17").Value = 110
$ws.Range("K17").Value = 24.545454545454
$ws.Range("L17").Value = 61.176470588235
$ws.Range("M17").Value = 144.642857142857
$ws.Range("N17").Value = -12.179487179487

$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 44.444444444444
$ws.Range("I18").Value = 79
$ws.Range("J18").Value = 116
$ws.Range("K18").Value = -31.896551724137
$ws.Range("L18").Value = -10.227272727272
$ws.Range("M18").Value = -45.517241379310
$ws.Range("N18").Value = -89.367429340511

$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -23.076923076923
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 50
$ws.Range("H19").Value = 2
$ws.Range("I19").Value = 321
$ws.Range("J19").Value = 336
$ws.Range("K19").Value = -4.464285714285
$ws.Range("L19").Value = -11.080332409972
$ws.Range("M19").Value = 45.909090909090
$ws.Range("N19").Value = -17.480719794344

$ws.Range("C20").Value = 6
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 63.636363636363
$ws.Range("I20").Value = 94
$ws.Range("J20").Value = 71
$ws.Range("K20").Value = 32.394366197183
$ws.Range("L20").Value = 36.231884057971
$ws.Range("M20").Value = 8.045977011494
$ws.Range("N20").Value = -89.782608695652

$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 4.545454545454
$ws.Range("F21").Value = 124
$ws.Range("G21").Value = 94
$ws.Range("H21").Value = 31.914893617021
$ws.Range("I21").Value = 712
$ws.Range("J21").Value = 714
$ws.Range("K21").Value = -0.280112044817
$ws.Range("L21").Value = 5.014749262536
$ws.Range("M21").Value = 18.469217970049
$ws.Range("N21").Value = -72.927756653992

$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 133.333333333333
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = 114.285714285714
$ws.Range("M22").Value = -16.666666666666

$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = -23.076923076923
$ws.Range("F24").Value = 105
$ws.Range("G24").Value = 147
$ws.Range("H24").Value = -28.571428571428
$ws.Range("I24").Value = 886
$ws.Range("J24").Value = 936
$ws.Range("K24").Value = -5.341880341880
$ws.Range("L24").Value = 1.026225769669
$ws.Range("M24").Value = 86.134453781512

$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = -22.222222222222
$ws.Range("F25").Value = 74
$ws.Range("G25").Value = 88
$ws.Range("H25").Value = -15.909090909090
$ws.Range("I25").Value = 510
$ws.Range("J25").Value = 545
$ws.Range("K25").Value = -6.422018348623
$ws.Range("L25").Value = 1.391650099403

$ws.Range("D26").Value = 9
$ws.Range("E26").Value = 55.555555555555
$ws.Range("F26").Value = 59
$ws.Range("G26").Value = 49
$ws.Range("H26").Value = 20.408163265306
$ws.Range("I26").Value = 336
$ws.Range("J26").Value = 309
$ws.Range("K26").Value = 8.737864077669
$ws.Range("L26").Value = 50
$ws.Range("M26").Value = 35.483870967741

$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 300
$ws.Range("I27").Value = 12
$ws.Range("K27").Value = 9.090909090909
$ws.Range("L27").Value = -7.692307692307

$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 7
$ws.Range("H28").Value = -22.222222222222
$ws.Range("J28").Value = 28
$ws.Range("K28").Value = -10.714285714285

$ws.Range("L29").Value = 0
$ws.Range("L30").Value = -50
$ws.Range("F31").Value = 5
$ws.Range("L31").Value = 100
